$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 467 (shifts rows 467:501 down to 468:502)
$ws.Rows.Item(467).Insert()

# Populate the new row 467 with the new record's data
$ws.Cells.Item(467, 1).Value = 4
$ws.Cells.Item(467, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(467, 3).Value = "Los Lagos"
$ws.Cells.Item(467, 4).Value = 45013
$ws.Cells.Item(467, 5).Value = 10
$ws.Cells.Item(467, 6).Value = 100112023
$ws.Cells.Item(467, 7).Value = "Brócoli"
$ws.Cells.Item(467, 8).Value = "Sin especificar"
$ws.Cells.Item(467, 9).Value = "Primera"
$ws.Cells.Item(467, 10).Value = 1200
$ws.Cells.Item(467, 11).Value = 1600
$ws.Cells.Item(467, 12).Value = 1600
$ws.Cells.Item(467, 13).Value = 1600
$ws.Cells.Item(467, 14).Value = "$/unidad"
$ws.Cells.Item(467, 15).Value = "Región Metropolitana"
$ws.Cells.Item(467, 16).Value = 1600
$ws.Cells.Item(467, 17).Value = 1
$ws.Cells.Item(467, 18).Value = "Hortaliza"
